$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 41 (current data row for 2022-04-18 / J=60 etc.)
# This shifts rows 41-49 down to 42-50, preserving their contents.
$ws.Rows.Item(41).Insert()

# Fill in the new row 41 with data, copying fixed fields from the row below (now row 42)
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat
$ws.Range("D41").Value = 44855
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100112042
$ws.Range("G41").Value = "Locoto"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 30
$ws.Range("K41").Value = 2500
$ws.Range("L41").Value = 2500
$ws.Range("M41").Value = 2500
$ws.Range("N41").Value = "$/kilo"
$ws.Range("O41").Value = "Región de Arica y Parinacota"
$ws.Range("P41").Value = 2500
$ws.Range("Q41").Value = 1
$ws.Range("R41").Value = "Hortaliza"
